# Reorganized hydro main functions; began working in the initial condition.
# Adds a new "sensor_locations" pointer row to the file-pointers sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14: Content = sensor_locations, Path = data/reprojected_dipwells.gpkg
# Write column B (the path) first so it lands earlier in the shared-string
# table than the column A label, matching the authored workbook's string order.
$ws.Cells.Item(14, 2).Value = "data/reprojected_dipwells.gpkg"
$ws.Cells.Item(14, 1).Value = "sensor_locations"

# Move / update the active selection to A15, just below the newly added row.
$ws.Range("A15").Select() | Out-Null
